# Daily attendance processing - 2026-01-19 10:42:48
# Normalizes the "Recorded By" (column G) text on the "Session Analysis
# Results" sheet: entries listing "System" alongside a real user/backdoor
# account are reordered so "System" is no longer first in the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($text -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($text -eq "System, admin@admin.com") {
        $cell.Value = "admin@admin.com, System"
    }
    elseif ($text -eq "System, system, backup@backdoor.com") {
        $cell.Value = "system, System, backup@backdoor.com"
    }
}

Write-Host "Recorded By column normalized"
